# Commit: "loop through eims mods per-cruise then assemble"
#
# Adds five new IODE quality-flag category rows to the CategoricalVariables
# sheet, and moves the "active sheet" selection from Keywords back to
# CategoricalVariables (which becomes the active tab / has the active
# selection in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CategoricalVariables")

# Make this the active sheet (also flips tabSelected off of whichever sheet
# was previously active, i.e. Keywords).
$ws.Activate()

# New attribute + its 5 valid codes/definitions (IODE quality flags).
$ws.Range("A4").Value = "iode_quality_flag"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "good"

$ws.Range("A5").Value = "iode_quality_flag"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "quality not evaluated, not available or unknown"

$ws.Range("A6").Value = "iode_quality_flag"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "questionable/suspect"

$ws.Range("A7").Value = "iode_quality_flag"
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = "bad"

$ws.Range("A8").Value = "iode_quality_flag"
$ws.Range("B8").Value = 9
$ws.Range("C8").Value = "missing data"

# Leave the selection on the newly added block, matching the saved view.
$ws.Range("A4:C8").Select()
